$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 12000
$ws.Range("J16").Value = 12000
$ws.Range("L16").Value = 12000
$ws.Range("N16").Value = -12460

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6480.654
$ws.Range("I51").Value = 10161.615
$ws.Range("J51").Value = 2799.6924
$ws.Range("K51").Value = 10161.615
$ws.Range("L51").Value = 2799.6924
$ws.Range("M51").Value = -9677.615
$ws.Range("N51").Value = -3767.6924

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 884.5714
$ws.Range("I92").Value = 906.64
$ws.Range("J92").Value = 829.4
$ws.Range("K92").Value = 906.64
$ws.Range("L92").Value = 829.4
$ws.Range("M92").Value = 341.36
$ws.Range("N92").Value = -3325.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1187.6
$ws.Range("I121").Value = 1333.3334
$ws.Range("J121").Value = 1125.1428
$ws.Range("K121").Value = 4000.0002
$ws.Range("L121").Value = 3375.4284
$ws.Range("M121").Value = -2253.0002
$ws.Range("N121").Value = -6869.428400000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3977.5625
$ws.Range("I131").Value = 528.5454999999999
$ws.Range("J131").Value = 5002.946
$ws.Range("K131").Value = 1585.6365
$ws.Range("L131").Value = 15008.838
$ws.Range("M131").Value = 3454.3635
$ws.Range("N131").Value = -25088.838

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8628138
$ws.Range("I132").Value = 8936263
$ws.Range("J132").Value = 650
$ws.Range("K132").Value = 26808789
$ws.Range("L132").Value = 1950
$ws.Range("M132").Value = -26806259
$ws.Range("N132").Value = -7010

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2349.9375
$ws.Range("I137").Value = 1300
$ws.Range("J137").Value = 3699.8572
$ws.Range("K137").Value = 3900
$ws.Range("L137").Value = 11099.5716
$ws.Range("M137").Value = -1350
$ws.Range("N137").Value = -16199.5716

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1521.2593
$ws.Range("I141").Value = 1356.0392
$ws.Range("J141").Value = 4330
$ws.Range("K141").Value = 4068.1176
$ws.Range("L141").Value = 12990
$ws.Range("M141").Value = 1111.8824
$ws.Range("N141").Value = -23350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29963.033
$ws.Range("I32").Value = 5548.107
$ws.Range("J32").Value = 257835.67
$ws.Range("K32").Value = 5548.107
$ws.Range("L32").Value = 257835.67
$ws.Range("M32").Value = -5261.107
$ws.Range("N32").Value = -258409.67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1446.8032
$ws.Range("I61").Value = 891.54285
$ws.Range("J61").Value = 2194.2693
$ws.Range("K61").Value = 891.54285
$ws.Range("L61").Value = 2194.2693
$ws.Range("M61").Value = -679.54285
$ws.Range("N61").Value = -2618.2693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 825.5714
$ws.Range("I74").Value = 800
$ws.Range("J74").Value = 829.8333
$ws.Range("K74").Value = 800
$ws.Range("L74").Value = 829.8333
$ws.Range("M74").Value = 74
$ws.Range("N74").Value = -2577.8333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 825.5714
$ws.Range("I77").Value = 800
$ws.Range("J77").Value = 829.8333
$ws.Range("K77").Value = 4000
$ws.Range("L77").Value = 4149.1665
$ws.Range("M77").Value = 368
$ws.Range("N77").Value = -12885.1665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1446.8032
$ws.Range("I136").Value = 891.54285
$ws.Range("J136").Value = 2194.2693
$ws.Range("K136").Value = 2674.62855
$ws.Range("L136").Value = 6582.8079
$ws.Range("M136").Value = -124.6285500000004
$ws.Range("N136").Value = -11682.8079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1984.7018
$ws.Range("I134").Value = 1794.4286
$ws.Range("K134").Value = 5383.2858
$ws.Range("M134").Value = -2848.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 19500
$ws.Range("I103").Value = 19500
$ws.Range("K103").Value = 19500
$ws.Range("M103").Value = -18328

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 21335.5
$ws.Range("J106").Value = 21335.5
$ws.Range("L106").Value = 21335.5
$ws.Range("N106").Value = -23859.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 7616.2666
$ws.Range("I107").Value = 15334.143
$ws.Range("J107").Value = 863.125
$ws.Range("K107").Value = 15334.143
$ws.Range("L107").Value = 863.125
$ws.Range("M107").Value = -13414.143
$ws.Range("N107").Value = -4703.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3148.8
$ws.Range("I132").Value = 2979.76
$ws.Range("K132").Value = 8939.280000000001
$ws.Range("M132").Value = -6409.280000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1387.1578
$ws.Range("I134").Value = 1371.8334
$ws.Range("J134").Value = 1413.4286
$ws.Range("K134").Value = 4115.5002
$ws.Range("L134").Value = 4240.2858
$ws.Range("M134").Value = -1580.5002
$ws.Range("N134").Value = -9310.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 125403.625
$ws.Range("I98").Value = 488.33334
$ws.Range("J98").Value = 500149.5
$ws.Range("K98").Value = 1465.00002
$ws.Range("L98").Value = 1500448.5
$ws.Range("M98").Value = 32.99998000000005
$ws.Range("N98").Value = -1503444.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1011.9722
$ws.Range("J131").Value = 1006.7656
$ws.Range("L131").Value = 3020.2968
$ws.Range("N131").Value = -13100.2968

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2076.574
$ws.Range("I132").Value = 2185.2666
$ws.Range("J132").Value = 1533.1111
$ws.Range("K132").Value = 6555.7998
$ws.Range("L132").Value = 4599.3333
$ws.Range("M132").Value = -4025.7998
$ws.Range("N132").Value = -9659.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2542.75
$ws.Range("J122").Value = 2446.6667
$ws.Range("L122").Value = 7340.000100000001
$ws.Range("N122").Value = -12240.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 28500
$ws.Range("J127").Value = 28500
$ws.Range("L127").Value = 28500
$ws.Range("N127").Value = -38420

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1029.7949
$ws.Range("I136").Value = 990.6111
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 2971.8333
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -421.8332999999998
$ws.Range("N136").Value = -9600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2145.28
$ws.Range("I122").Value = 1751.9286
$ws.Range("J122").Value = 2645.9092
$ws.Range("K122").Value = 5255.7858
$ws.Range("L122").Value = 7937.7276
$ws.Range("M122").Value = -2805.7858
$ws.Range("N122").Value = -12837.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2335.5
$ws.Range("I126").Value = 2245
$ws.Range("K126").Value = 6735
$ws.Range("M126").Value = -4265

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1628.9546
$ws.Range("I132").Value = 1682.644
$ws.Range("J132").Value = 1176.4286
$ws.Range("K132").Value = 5047.932
$ws.Range("L132").Value = 3529.2858
$ws.Range("M132").Value = -2517.932
$ws.Range("N132").Value = -8589.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 606.8627300000001
$ws.Range("I136").Value = 371.79486
$ws.Range("J136").Value = 1370.8334
$ws.Range("K136").Value = 1115.38458
$ws.Range("L136").Value = 4112.5002
$ws.Range("M136").Value = 1434.61542
$ws.Range("N136").Value = -9212.5002
